# Updated cryptos list on Sat Oct  5 21:51:54 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.885.52'
$ws.Range("E2").Value = '  -0.72%  '
$ws.Range("D3").Value = '2.407.41'
$ws.Range("E3").Value = '  -0.70%  '
$ws.Range("E4").Value = '  +0.01%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '561.72'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.95%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '141.97'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -1.25%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  -0.88%  '
$ws.Range("E9").Value = '  -0.91%  '
$ws.Range("E11").Value = '  -2.92%  '
$ws.Range("E12").Value = '  -0.84%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '25.46'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -3.09%  '
$ws.Range("D15").Value = '2.843.88'
$ws.Range("E15").Value = '  -0.61%  '
$ws.Range("D16").Value = '61.751.12'
$ws.Range("E16").Value = '  -0.70%  '
$ws.Range("D17").Value = '2.404.26'
$ws.Range("E17").Value = '  -0.77%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '11.21'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +1.12%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '321.00'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -1.09%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '6.81'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +1.14%  '
$ws.Range("E21").Value = '  -1.87%  '
$ws.Range("E22").Value = '  -0.21%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '65.64'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +1.22%  '
$ws.Range("E24").Value = '  -2.65%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '8.66'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -4.93%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '564.06'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -1.61%  '
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("D28").Value = '2.515.27'
$ws.Range("D29").Value = '0.0₃0928'
$ws.Range("E29").Value = '  -1.14%  '
$ws.Range("E30").Value = '  -2.73%  '
$ws.Range("E31").Value = '  -4.80%  '
$ws.Range("E32").Value = '  -0.81%  '
$ws.Range("E33").Value = '  +0.32%  '
$ws.Range("E34").Value = '  -4.10%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("E36").Value = '  -1.92%  '
$ws.Range("E39").Value = '  -1.55%  '
$ws.Range("E40").Value = '  -1.55%  '
$ws.Range("E41").Value = '  -5.29%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '147.93'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -2.20%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '2.23'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -4.25%  '
$ws.Range("E45").Value = '  -1.42%  '
$ws.Range("E46").Value = '  -2.32%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '19.78'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -2.88%  '
$ws.Range("E48").Value = '  -0.11%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.0917'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +0.12%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.0224'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -1.81%  '
$ws.Range("E51").Value = '  +0.38%  '
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '152.61'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +1.79%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '5.42'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -5.73%  '
